$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "Normal" cell style explicitly to every populated cell in the
# used range. This mirrors the new second cellXfs entry that the edited
# workbook references for all its data cells (style index 1), while
# leaving genuinely empty cells untouched (no phantom <c> entries).
$used = $ws.Range("A1:S7").SpecialCells(2)
foreach ($area in $used.Areas) {
    $area.Style = "Normal"
}

# Update row 7 (G7:I7, L7) from the old "production of bike" / DK flow to
# the new "carbon fibre production" / DE flow, reusing the existing
# shared-string text so the workbook's string table stays deduplicated.
$ws.Range("G7").Value = "carbon fibre production"
$ws.Range("H7").Value = "carbon fibre"
$ws.Range("I7").Value = "DE"
$ws.Range("L7").Value = "('bike_production_example', '7bde0c388a1b401c95dccd0a3429bd0c')"

# Move the active selection to G7:L7 with G7 as the active cell.
$ws.Range("G7:L7").Select() | Out-Null
